# Reorder/refresh the "Periodo Mora" (column E) and "Valor Mora" (column F)
# data rows 16-28 on Hoja1: the account-statement periods are now listed in
# reverse chronological order (most recent period first), and each period's
# "Valor Mora" value travels together with its period label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow  = 28

# Capture the current (pre-edit) Periodo Mora / Valor Mora pairs, top to bottom.
$periods = @()
$valores = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $valores += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse order so row 16 now holds the most recent
# period (previously on row 28) and row 28 holds the oldest (previously on
# row 16), carrying each row's own "Valor Mora" along with it.
$n = $periods.Count
for ($i = 0; $i -lt $n; $i++) {
    $targetRow = $firstRow + $i
    $srcIndex  = $n - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value2 = $periods[$srcIndex]
    $ws.Cells.Item($targetRow, 6).Value2 = $valores[$srcIndex]
}
